$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "32"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "16"
$ws2.PageSetup.LeftMargin = 56.7
$ws2.PageSetup.RightMargin = 56.7
$ws2.PageSetup.TopMargin = 75.8
$ws2.PageSetup.BottomMargin = 75.8
$ws2.PageSetup.HeaderMargin = 56.7
$ws2.PageSetup.FooterMargin = 56.7
$ws2.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws2.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
$ws2.PageSetup.Orientation = 1
$ws2.PageSetup.PaperSize = 1
$ws2.PageSetup.Zoom = 100
Write-Output "done"
